$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Cells.Item(33, 8).Value = 745.5625
$ws.Cells.Item(33, 9).Value = 447.92307
$ws.Cells.Item(33, 10).Value = 2035.3334
$ws.Cells.Item(33, 11).Value = 447.92307
$ws.Cells.Item(33, 12).Value = 2035.3334
$ws.Cells.Item(33, 13).Value = -218.92307
$ws.Cells.Item(33, 14).Value = -2493.3334

$ws.Cells.Item(70, 8).Value = 2910.1667
$ws.Cells.Item(70, 9).Value = 1571.8889
$ws.Cells.Item(70, 10).Value = 3483.7144
$ws.Cells.Item(70, 11).Value = 4715.6667
$ws.Cells.Item(70, 12).Value = 10451.1432
$ws.Cells.Item(70, 13).Value = -4445.6667
$ws.Cells.Item(70, 14).Value = -10991.1432

$ws.Cells.Item(73, 8).Value = 2910.1667
$ws.Cells.Item(73, 9).Value = 1571.8889
$ws.Cells.Item(73, 10).Value = 3483.7144
$ws.Cells.Item(73, 11).Value = 4715.6667
$ws.Cells.Item(73, 12).Value = 10451.1432
$ws.Cells.Item(73, 13).Value = -3779.6667
$ws.Cells.Item(73, 14).Value = -12323.1432

$ws.Cells.Item(76, 8).Value = 2000
$ws.Cells.Item(76, 9).Value = 2000
$ws.Cells.Item(76, 10).Value = 2000
$ws.Cells.Item(76, 11).Value = 2000
$ws.Cells.Item(76, 12).Value = 2000
$ws.Cells.Item(76, 13).Value = -1685
$ws.Cells.Item(76, 14).Value = -2630

$ws.Cells.Item(79, 8).Value = 2000
$ws.Cells.Item(79, 9).Value = 2000
$ws.Cells.Item(79, 10).Value = 2000
$ws.Cells.Item(79, 11).Value = 2000
$ws.Cells.Item(79, 12).Value = 2000
$ws.Cells.Item(79, 13).Value = -908
$ws.Cells.Item(79, 14).Value = -4184

$ws.Cells.Item(132, 8).Value = 8002555.5
$ws.Cells.Item(132, 9).Value = 8335579
$ws.Cells.Item(132, 11).Value = 25006737
$ws.Cells.Item(132, 13).Value = -25004207

$ws.Cells.Item(138, 8).Value = 3832.1272
$ws.Cells.Item(138, 9).Value = 1777.88
$ws.Cells.Item(138, 10).Value = 5544
$ws.Cells.Item(138, 11).Value = 5333.64
$ws.Cells.Item(138, 12).Value = 16632
$ws.Cells.Item(138, 13).Value = -193.6400000000003
$ws.Cells.Item(138, 14).Value = -26912

$ws = $wb.Worksheets("ARM")
$ws.Cells.Item(2, 8).Value = 6252411.5
$ws.Cells.Item(2, 9).Value = 13890162
$ws.Cells.Item(2, 10).Value = 3343.3635
$ws.Cells.Item(2, 11).Value = 13890162
$ws.Cells.Item(2, 12).Value = 3343.3635
$ws.Cells.Item(2, 13).Value = -13890049
$ws.Cells.Item(2, 14).Value = -3569.3635

$ws.Cells.Item(74, 8).Value = 1075.375
$ws.Cells.Item(74, 9).Value = 1017.3333
$ws.Cells.Item(74, 10).Value = 1249.5
$ws.Cells.Item(74, 11).Value = 1017.3333
$ws.Cells.Item(74, 12).Value = 1249.5
$ws.Cells.Item(74, 13).Value = -143.3333
$ws.Cells.Item(74, 14).Value = -2997.5

$ws.Cells.Item(77, 8).Value = 1075.375
$ws.Cells.Item(77, 9).Value = 1017.3333
$ws.Cells.Item(77, 10).Value = 1249.5
$ws.Cells.Item(77, 11).Value = 5086.6665
$ws.Cells.Item(77, 12).Value = 6247.5
$ws.Cells.Item(77, 13).Value = -718.6665000000003
$ws.Cells.Item(77, 14).Value = -14983.5

$ws.Cells.Item(116, 8).Value = 6252411.5
$ws.Cells.Item(116, 9).Value = 13890162
$ws.Cells.Item(116, 10).Value = 3343.3635
$ws.Cells.Item(116, 11).Value = 13890162
$ws.Cells.Item(116, 12).Value = 3343.3635
$ws.Cells.Item(116, 13).Value = -13887868
$ws.Cells.Item(116, 14).Value = -7931.363499999999

$ws.Cells.Item(141, 8).Value = 40000
$ws.Cells.Item(141, 10).Value = 40000
$ws.Cells.Item(141, 12).Value = 40000
$ws.Cells.Item(141, 14).Value = -50360

$ws = $wb.Worksheets("BSM")
$ws.Cells.Item(3, 8).Value = 6252411.5
$ws.Cells.Item(3, 9).Value = 13890162
$ws.Cells.Item(3, 10).Value = 3343.3635
$ws.Cells.Item(3, 11).Value = 13890162
$ws.Cells.Item(3, 12).Value = 3343.3635
$ws.Cells.Item(3, 13).Value = -13890048
$ws.Cells.Item(3, 14).Value = -3571.3635

$ws.Cells.Item(86, 8).Value = 2600
$ws.Cells.Item(86, 9).Value = 1442.8572
$ws.Cells.Item(86, 10).Value = 3410
$ws.Cells.Item(86, 11).Value = 1442.8572
$ws.Cells.Item(86, 12).Value = 3410
$ws.Cells.Item(86, 13).Value = -319.8571999999999
$ws.Cells.Item(86, 14).Value = -5656

$ws.Cells.Item(89, 8).Value = 2600
$ws.Cells.Item(89, 9).Value = 1442.8572
$ws.Cells.Item(89, 10).Value = 3410
$ws.Cells.Item(89, 11).Value = 7214.286
$ws.Cells.Item(89, 12).Value = 17050
$ws.Cells.Item(89, 13).Value = -1598.286
$ws.Cells.Item(89, 14).Value = -28282

$ws.Cells.Item(134, 8).Value = 2753.16
$ws.Cells.Item(134, 9).Value = 1942.8823
$ws.Cells.Item(134, 10).Value = 4475
$ws.Cells.Item(134, 11).Value = 5828.6469
$ws.Cells.Item(134, 12).Value = 13425
$ws.Cells.Item(134, 13).Value = -3293.6469
$ws.Cells.Item(134, 14).Value = -18495

$ws = $wb.Worksheets("CRP")
$ws.Cells.Item(31, 8).Value = 3156.558
$ws.Cells.Item(31, 9).Value = 1746.8387
$ws.Cells.Item(31, 10).Value = 6798.3335
$ws.Cells.Item(31, 11).Value = 1746.8387
$ws.Cells.Item(31, 12).Value = 6798.3335
$ws.Cells.Item(31, 13).Value = -1451.8387
$ws.Cells.Item(31, 14).Value = -7388.3335

$ws.Cells.Item(34, 8).Value = 3156.558
$ws.Cells.Item(34, 9).Value = 1746.8387
$ws.Cells.Item(34, 10).Value = 6798.3335
$ws.Cells.Item(34, 11).Value = 1746.8387
$ws.Cells.Item(34, 12).Value = 6798.3335
$ws.Cells.Item(34, 13).Value = -1544.8387
$ws.Cells.Item(34, 14).Value = -7202.3335

$ws.Cells.Item(58, 8).Value = 19234074
$ws.Cells.Item(58, 9).Value = 2333.2666
$ws.Cells.Item(58, 10).Value = 45459172
$ws.Cells.Item(58, 11).Value = 2333.2666
$ws.Cells.Item(58, 12).Value = 45459172
$ws.Cells.Item(58, 13).Value = -2130.2666
$ws.Cells.Item(58, 14).Value = -45459578

$ws.Cells.Item(86, 8).Value = 5543.357
$ws.Cells.Item(86, 9).Value = 3511.889
$ws.Cells.Item(86, 11).Value = 3511.889
$ws.Cells.Item(86, 13).Value = -2388.889

$ws.Cells.Item(89, 8).Value = 5543.357
$ws.Cells.Item(89, 9).Value = 3511.889
$ws.Cells.Item(89, 11).Value = 17559.445
$ws.Cells.Item(89, 13).Value = -11943.445

$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 19234074
$ws.Cells.Item(136, 9).Value = 2333.2666
$ws.Cells.Item(136, 10).Value = 45459172
$ws.Cells.Item(136, 11).Value = 6999.7998
$ws.Cells.Item(136, 12).Value = 136377516
$ws.Cells.Item(136, 13).Value = -4449.7998
$ws.Cells.Item(136, 14).Value = -136382616

$ws = $wb.Worksheets("CUL")
$ws.Cells.Item(34, 8).Value = 11916.5
$ws.Cells.Item(34, 9).Value = 185
$ws.Cells.Item(34, 10).Value = 13220
$ws.Cells.Item(34, 11).Value = 555
$ws.Cells.Item(34, 12).Value = 39660
$ws.Cells.Item(34, 13).Value = -471
$ws.Cells.Item(34, 14).Value = -39828

$ws.Cells.Item(94, 8).Value = 3668.1516
$ws.Cells.Item(94, 10).Value = 3720.2812
$ws.Cells.Item(94, 12).Value = 11160.8436
$ws.Cells.Item(94, 14).Value = -12512.8436

$ws.Cells.Item(131, 8).Value = 1066.326
$ws.Cells.Item(131, 9).Value = 848.1429000000001
$ws.Cells.Item(131, 10).Value = 1405.7222
$ws.Cells.Item(131, 11).Value = 2544.4287
$ws.Cells.Item(131, 12).Value = 4217.1666
$ws.Cells.Item(131, 13).Value = 2495.5713
$ws.Cells.Item(131, 14).Value = -14297.1666

$ws = $wb.Worksheets("GSM")
$ws.Cells.Item(70, 8).Value = 5759.3
$ws.Cells.Item(70, 9).Value = 6325.933
$ws.Cells.Item(70, 10).Value = 4059.4
$ws.Cells.Item(70, 11).Value = 6325.933
$ws.Cells.Item(70, 12).Value = 4059.4
$ws.Cells.Item(70, 13).Value = -6055.933
$ws.Cells.Item(70, 14).Value = -4599.4

$ws.Cells.Item(73, 8).Value = 5759.3
$ws.Cells.Item(73, 9).Value = 6325.933
$ws.Cells.Item(73, 10).Value = 4059.4
$ws.Cells.Item(73, 11).Value = 6325.933
$ws.Cells.Item(73, 12).Value = 4059.4
$ws.Cells.Item(73, 13).Value = -5389.933
$ws.Cells.Item(73, 14).Value = -5931.4

$ws.Cells.Item(122, 8).Value = 4553.1665
$ws.Cells.Item(122, 9).Value = 5914.143
$ws.Cells.Item(122, 11).Value = 17742.429
$ws.Cells.Item(122, 13).Value = -15292.429

$ws = $wb.Worksheets("LTW")
$ws.Cells.Item(14, 8).Value = 324879
$ws.Cells.Item(14, 10).Value = 21255
$ws.Cells.Item(14, 12).Value = 21255
$ws.Cells.Item(14, 14).Value = -21599

$ws.Cells.Item(35, 8).Value = 9277.5
$ws.Cells.Item(35, 9).Value = 2370
$ws.Cells.Item(35, 11).Value = 2370
$ws.Cells.Item(35, 13).Value = -2034

$ws.Cells.Item(46, 8).Value = 2162.8667
$ws.Cells.Item(46, 9).Value = 673.3333
$ws.Cells.Item(46, 10).Value = 3155.889
$ws.Cells.Item(46, 11).Value = 673.3333
$ws.Cells.Item(46, 12).Value = 3155.889
$ws.Cells.Item(46, 13).Value = -485.3333
$ws.Cells.Item(46, 14).Value = -3531.889

$ws.Cells.Item(82, 8).Value = 2393.074
$ws.Cells.Item(82, 9).Value = 1945.7142
$ws.Cells.Item(82, 10).Value = 2874.8462
$ws.Cells.Item(82, 11).Value = 1945.7142
$ws.Cells.Item(82, 12).Value = 2874.8462
$ws.Cells.Item(82, 13).Value = -1584.7142
$ws.Cells.Item(82, 14).Value = -3596.8462

$ws.Cells.Item(85, 8).Value = 2393.074
$ws.Cells.Item(85, 9).Value = 1945.7142
$ws.Cells.Item(85, 10).Value = 2874.8462
$ws.Cells.Item(85, 11).Value = 1945.7142
$ws.Cells.Item(85, 12).Value = 2874.8462
$ws.Cells.Item(85, 13).Value = -697.7141999999999
$ws.Cells.Item(85, 14).Value = -5370.8462

$ws.Cells.Item(93, 8).Value = 1573.258
$ws.Cells.Item(93, 9).Value = 1095.85
$ws.Cells.Item(93, 10).Value = 2441.2727
$ws.Cells.Item(93, 11).Value = 1095.85
$ws.Cells.Item(93, 12).Value = 2441.2727
$ws.Cells.Item(93, 13).Value = 152.1500000000001
$ws.Cells.Item(93, 14).Value = -4937.2727

$ws.Cells.Item(132, 8).Value = 3569.077
$ws.Cells.Item(132, 9).Value = 6000
$ws.Cells.Item(132, 10).Value = 3366.5
$ws.Cells.Item(132, 11).Value = 18000
$ws.Cells.Item(132, 12).Value = 10099.5
$ws.Cells.Item(132, 13).Value = -15470
$ws.Cells.Item(132, 14).Value = -15159.5

$ws = $wb.Worksheets("WVR")
$ws.Cells.Item(62, 8).Value = 3195
$ws.Cells.Item(62, 9).Value = 3390
$ws.Cells.Item(62, 11).Value = 3390
$ws.Cells.Item(62, 13).Value = -2766

$ws.Cells.Item(65, 8).Value = 3195
$ws.Cells.Item(65, 9).Value = 3390
$ws.Cells.Item(65, 11).Value = 16950
$ws.Cells.Item(65, 13).Value = -13830

$ws.Cells.Item(69, 8).Value = 69125
$ws.Cells.Item(69, 10).Value = 69125
$ws.Cells.Item(69, 12).Value = 69125
$ws.Cells.Item(69, 14).Value = -70623

$ws.Cells.Item(72, 8).Value = 69125
$ws.Cells.Item(72, 10).Value = 69125
$ws.Cells.Item(72, 12).Value = 207375
$ws.Cells.Item(72, 14).Value = -214863

$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

Write-Host "Applied 37 row updates across 8 sheets"
